$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.557.76"
$ws.Range("E2").Value = "  -1.31%  "
$ws.Range("D3").Value = "3.011.73"
$ws.Range("E3").Value = "  -1.64%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'585.08"
$ws.Range("E5").Value = "  -1.29%  "
$ws.Range("D6").Value = "'146.54"
$ws.Range("E6").Value = "  -4.96%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.527"
$ws.Range("E8").Value = "  -2.55%  "
$ws.Range("D9").Value = "3.011.02"
$ws.Range("E9").Value = "  -1.71%  "
$ws.Range("E10").Value = "  -4.10%  "
$ws.Range("D11").Value = "'5.82"
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("D12").Value = "'0.460"
$ws.Range("E12").Value = "  +1.92%  "
$ws.Range("E13").Value = "  -3.08%  "
$ws.Range("D14").Value = "'34.72"
$ws.Range("E14").Value = "  -5.69%  "
$ws.Range("E15").Value = "  +1.99%  "
$ws.Range("D16").Value = "3.505.36"
$ws.Range("E16").Value = "  -1.90%  "
$ws.Range("D17").Value = "'7.10"
$ws.Range("E17").Value = "  -1.08%  "
$ws.Range("D18").Value = "62.480.74"
$ws.Range("E18").Value = "  -1.42%  "
$ws.Range("D19").Value = "3.012.02"
$ws.Range("E19").Value = "  -1.77%  "
$ws.Range("D20").Value = "'459.15"
$ws.Range("E20").Value = "  -6.07%  "
$ws.Range("D21").Value = "'13.99"
$ws.Range("E21").Value = "  -2.79%  "
$ws.Range("D22").Value = "'0.689"
$ws.Range("E22").Value = "  -2.60%  "
$ws.Range("E23").Value = "  -1.55%  "
$ws.Range("D24").Value = "'81.66"
$ws.Range("E24").Value = "  -0.52%  "
$ws.Range("E25").Value = "  -3.90%  "
$ws.Range("E26").Value = "  -8.75%  "
$ws.Range("E27").Value = "  -6.32%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.20%  "
$ws.Range("E30").Value = "  -2.62%  "
$ws.Range("D31").Value = "'7.06"
$ws.Range("E31").Value = "  -4.32%  "
$ws.Range("E32").Value = "  -5.74%  "
$ws.Range("D33").Value = "'28.04"
$ws.Range("E33").Value = "  +2.34%  "
$ws.Range("E34").Value = "  -2.18%  "
$ws.Range("D35").Value = "0.0₃0812"
$ws.Range("E35").Value = "  -1.31%  "
$ws.Range("E36").Value = "  -3.57%  "
$ws.Range("D37").Value = "'5.78"
$ws.Range("E37").Value = "  -3.54%  "
$ws.Range("E38").Value = "  -5.24%  "
$ws.Range("D39").Value = "'50.34"
$ws.Range("E39").Value = "  -0.53%  "
$ws.Range("D40").Value = "'9.13"
$ws.Range("E40").Value = "  -1.64%  "
$ws.Range("D41").Value = "'2.91"
$ws.Range("E41").Value = "  -12.81%  "
$ws.Range("D42").Value = "'0.119"
$ws.Range("E42").Value = "  +4.71%  "
$ws.Range("D43").Value = "'389.97"
$ws.Range("E43").Value = "  -11.22%  "
$ws.Range("E44").Value = "  -1.55%  "
$ws.Range("E45").Value = "  -7.63%  "
$ws.Range("D46").Value = "2.736.19"
$ws.Range("E46").Value = "  -4.07%  "
$ws.Range("D47").Value = "'37.52"
$ws.Range("E47").Value = "  -3.48%  "
$ws.Range("D48").Value = "'129.49"
$ws.Range("E48").Value = "  -0.60%  "
$ws.Range("D50").Value = "'2.20"
$ws.Range("E50").Value = "  -1.52%  "
$ws.Range("E51").Value = "  -0.81%  "
